# Apply updated dSF (column F) values to Sheet1, per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F2"  = 3
    "F4"  = -1
    "F5"  = -1
    "F9"  = -2
    "F11" = -1
    "F19" = -1
    "F23" = -6
    "F24" = -3
    "F27" = -9
    "F30" = -3
    "F31" = 7
    "F32" = -2
    "F36" = -5
    "F40" = -2
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
